$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) AddCommodities sheet: append six new "Invalid_AddCommodity_TC001" rows
# ---------------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item("AddCommodities")

# Column A (Automation Test ID) uses the same style as the existing A2/A3
# cells (vertical-center alignment).
$wsAdd.Range("A4:A9").VerticalAlignment = -4108

# Row 4
$wsAdd.Range("A4").Value = "Invalid_AddCommodity_TC001"
$wsAdd.Range("B4").Value = "NA"
$wsAdd.Range("C4").Value = "NA"
$wsAdd.Range("D4").Value = "NA"
$wsAdd.Range("E4").Value = "NA"
$wsAdd.Range("F4").Value = "Add"
$wsAdd.Range("G4").Value = "Unable to add commodity"

# Row 5
$wsAdd.Range("A5").Value = "Invalid_AddCommodity_TC001(2)"
$wsAdd.Range("B5").Value = "Com"
$wsAdd.Range("C5").Value = "NA"
$wsAdd.Range("D5").Value = "NA"
$wsAdd.Range("E5").Value = "NA"
$wsAdd.Range("F5").Value = "Add"
$wsAdd.Range("G5").Value = "the upper limit is required"

# Row 6
$wsAdd.Range("A6").Value = "Invalid_AddCommodity_TC001(3)"
$wsAdd.Range("B6").Value = "Com"
$wsAdd.Range("C6").Value = 5000
$wsAdd.Range("C6").NumberFormat = "@"
$wsAdd.Range("D6").Value = "NA"
$wsAdd.Range("E6").Value = "NA"
$wsAdd.Range("F6").Value = "Add"
$wsAdd.Range("G6").Value = "the lower limit is required"

# Row 7
$wsAdd.Range("A7").Value = "Invalid_AddCommodity_TC001(4)"
$wsAdd.Range("B7").Value = "Com"
$wsAdd.Range("C7").NumberFormat = "@"
$wsAdd.Range("C7").Value = "NA"
$wsAdd.Range("D7").Value = 6000
$wsAdd.Range("D7").NumberFormat = "@"
$wsAdd.Range("E7").Value = "NA"
$wsAdd.Range("F7").Value = "Add"
$wsAdd.Range("G7").Value = "the upper limit is required"

# Row 8
$wsAdd.Range("A8").Value = "Invalid_AddCommodity_TC001(5)"
$wsAdd.Range("B8").Value = "NA"
$wsAdd.Range("C8").NumberFormat = "@"
$wsAdd.Range("C8").Value = "7000"
$wsAdd.Range("D8").Value = 6000
$wsAdd.Range("D8").NumberFormat = "@"
$wsAdd.Range("E8").Value = "NA"
$wsAdd.Range("F8").Value = "Add"
$wsAdd.Range("G8").Value = "the commodity name is required"

# Row 9
$wsAdd.Range("A9").Value = "Invalid_AddCommodity_TC001(6)"
$wsAdd.Range("B9").Value = "Com"
$wsAdd.Range("C9").NumberFormat = "@"
$wsAdd.Range("C9").Value = "7000"
$wsAdd.Range("D9").NumberFormat = "@"
$wsAdd.Range("D9").Value = "8000"
$wsAdd.Range("E9").Value = "NA"
$wsAdd.Range("F9").Value = "Add"
$wsAdd.Range("G9").Value = "the upper limit must be greater than the lower limit"

# ---------------------------------------------------------------------------
# 2) CustomizeGrid sheet: append a matching summary row
# ---------------------------------------------------------------------------
$wsGrid = $wb.Worksheets.Item("CustomizeGrid")
$wsGrid.Range("A4").VerticalAlignment = -4108
$wsGrid.Range("A4").Value = "Invalid_AddCommodity_TC001"
$wsGrid.Range("B4").Value = "NA"
$wsGrid.Range("C4").Value = "ALL"
$wsGrid.Range("D4").Value = "Webtable customized successfully"

# Select A4 on this sheet (matches the authored selection), without leaving
# it as the active tab.
$wsGrid.Activate()
$wsGrid.Range("A4").Select()

# ---------------------------------------------------------------------------
# 3) CommoditiesGrid sheet: select entire row 2
# ---------------------------------------------------------------------------
$wsCG = $wb.Worksheets.Item("CommoditiesGrid")
$wsCG.Activate()
$wsCG.Rows(2).Select()

# ---------------------------------------------------------------------------
# 4) EditCommodities sheet: keep its own B2 selection (no longer active tab)
# ---------------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("EditCommodities")
$wsEdit.Activate()
$wsEdit.Range("B2").Select()

# ---------------------------------------------------------------------------
# 5) Finish on AddCommodities, selecting G9 - this becomes the active tab
# ---------------------------------------------------------------------------
$wsAdd.Activate()
$wsAdd.Range("G9").Select()
